# Remove the trailing "Ver no Jupiter ..." / copyright footer block that
# was stripped from the bottom of the document (site footer removed from
# the Jekyll build). Three paragraphs go away:
#   - the blank paragraph right after the "LOQ4086..." requisito line
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#     pages. Original theme under Creative Commons Attribution"
# The paragraph that used to hold "LOQ4086..." stays untouched, as does
# the blank paragraph that precedes the final page-break paragraph.

$d = $word.ActiveDocument

$target = "Ver no Jupiter Salvar em pdf Salvar em docx"
$found = $false

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $target) {
        $found = $true

        # Delete this paragraph ("Ver no Jupiter ...")
        $p.Range.Delete()

        # Delete the next paragraph (the copyright/footer line), now at
        # the same index $i since the previous delete shifted things up.
        $copyright = $d.Paragraphs.Item($i)
        $copyright.Range.Delete()

        # Delete the blank paragraph immediately preceding "Ver no
        # Jupiter ..." (index $i - 1, before the two deletions above).
        $blank = $d.Paragraphs.Item($i - 1)
        $blank.Range.Delete()

        break
    }
}

if (-not $found) {
    throw "Could not locate the 'Ver no Jupiter ...' paragraph to remove"
}

Write-Output "Paragraphs remaining: $($d.Paragraphs.Count)"
